$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.144.20"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "'1.901.47"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'252.89"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "'0.698"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'41.67"
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("D9").Value = "'0.355"
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("D10").Value = "'52.40"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "'0.0757"
$ws.Range("E11").Value = "  +4.93%  "
$ws.Range("D12").Value = "'0.0979"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "'13.19"
$ws.Range("E13").Value = "  +5.10%  "
$ws.Range("D14").Value = "'2.180.25"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("E15").Value = "  +3.66%  "
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").Value = "'1.910.31"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "'35.150.06"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'73.69"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "'0.0₃0841"
$ws.Range("E20").Value = "  +2.45%  "
$ws.Range("D21").Value = "'243.16"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "'13.04"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "'5.02"
$ws.Range("E23").Value = "  +4.84%  "
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("E25").Value = "  +4.74%  "
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").Value = "'169.04"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").Value = "'8.57"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "'4.128.38"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  +11.92%  "
$ws.Range("E33").Value = "  +4.37%  "
$ws.Range("D34").Value = "'0.0594"
$ws.Range("E34").Value = "  +4.54%  "
$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  +6.32%  "
$ws.Range("D36").Value = "'4.24"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").Value = "'0.848"
$ws.Range("E38").Value = "  -7.29%  "
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "'17.33"
$ws.Range("E40").Value = "  +5.49%  "
$ws.Range("D41").Value = "'98.28"
$ws.Range("E41").Value = "  +5.79%  "
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("D43").Value = "'0.0664"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").Value = "'1.305.11"
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Value = "'2.74"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").Value = "'6.58"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").Value = "'12.00"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("E51").Value = "  +6.92%  "
